# Legacy GSC export data refresh.
#
# The "Chart" sheet holds a daily Date/Not-indexed/Indexed/Impressions
# table that starts one day too early (2025-10-26). The upstream export
# dropped that oldest day, so every subsequent day's Not-indexed /
# Indexed / Impressions figures shift up by one row while the calendar
# date column itself is untouched - i.e. the fix is simply removing the
# now-stale first data row (row 2, right under the header).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows("2").Delete()
